$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - copy formatting (bold/border/centered) from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cells H2/H3 with numeric value 1 (plain, unstyled like F2/F3, G2/G3)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1

$excel.CutCopyMode = 0
